$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-09-18 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-09-19 Thursday", 2)

# Update the multiplication problems in the table, cell by cell (addressed
# by row/column so that values which collide between old/new text across
# different cells -- e.g. "92×15=" -- are never mismatched by a global
# find/replace).
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="94×22="},
    @{Row=1;  Col=2; Text="56×19="},
    @{Row=1;  Col=3; Text="25×72="},
    @{Row=1;  Col=4; Text="93×62="},
    @{Row=1;  Col=5; Text="15×99="},

    @{Row=5;  Col=1; Text="28×47="},
    @{Row=5;  Col=2; Text="65×66="},
    @{Row=5;  Col=3; Text="92×15="},
    @{Row=5;  Col=4; Text="46×67="},
    @{Row=5;  Col=5; Text="93×81="},

    @{Row=10; Col=1; Text="77×19="},
    @{Row=10; Col=2; Text="31×76="},
    @{Row=10; Col=3; Text="88×86="},
    @{Row=10; Col=4; Text="26×44="},
    @{Row=10; Col=5; Text="90×45="},

    @{Row=15; Col=1; Text="61×31="},
    @{Row=15; Col=2; Text="83×63="},
    @{Row=15; Col=3; Text="89×91="},
    @{Row=15; Col=4; Text="57×45="},
    @{Row=15; Col=5; Text="67×22="},

    @{Row=20; Col=1; Text="93×88="},
    @{Row=20; Col=2; Text="24×31="},
    @{Row=20; Col=3; Text="33×29="},
    @{Row=20; Col=4; Text="26×38="},
    @{Row=20; Col=5; Text="67×77="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
